# REFACTOR: Bring all of Dokeza up to date including the online version.
#
# The document was carrying a stale SharePoint "Document Information Panel"
# custom XML data store (content-type schema, SharePoint form template and
# document-management properties, plus their datastore item-property
# parts). None of it is surfaced in the document body - it is pure
# library/content-type plumbing left over from the old SharePoint list -
# so bringing the file up to date means stripping all of the custom XML
# parts from the package.
#
# Real Word COM automation for this is: walk Document.CustomXMLParts
# (optionally scoping with SelectByNamespace/SelectByID) and call
# .Delete on every part that isn't one Word manages itself.

$d = $word.ActiveDocument

# Namespaces used by the legacy Dokeza/SharePoint custom XML parts that
# must be removed.
$targetNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)

$parts = $d.CustomXMLParts
Write-Host "CustomXMLParts found:" $parts.Count

# Walk backwards since deleting shifts the collection's indices.
for ($i = $parts.Count; $i -ge 1; $i--) {
    $part = $parts.Item($i)

    $ns = $part.NamespaceURI
    if ($targetNamespaces -contains $ns) {
        Write-Host "Deleting custom XML part" $i "(" $ns ")"
        $part.Delete()
    }
}

# Belt-and-braces: some hosts only expose matching parts through
# SelectByNamespace rather than raw iteration, so sweep that path too in
# case any of the target parts were left behind above.
foreach ($ns in $targetNamespaces) {
    try {
        $scoped = $d.CustomXMLParts.SelectByNamespace($ns)
        for ($j = $scoped.Count; $j -ge 1; $j--) {
            $scoped.Item($j).Delete()
        }
    } catch {
        # Namespace not present / already removed - nothing to do.
    }
}

Write-Host "CustomXMLParts remaining:" $d.CustomXMLParts.Count

$d.Save()
